$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for row 4 (theta_se) and row 6 (lambda_se), replacing the old
# "(nan)" placeholders. The shared-strings table in the target file was
# built column-by-column (row4 then row6 for each column in turn), so we
# set the cells in that same interleaved order to reproduce it.
$ws.Range("B4").Value = "(0.33)"
$ws.Range("B6").Value = "(0.57)"

$ws.Range("C4").Value = "(1.22)"
$ws.Range("C6").Value = "(0.75)"

$ws.Range("D4").Value = "(0.4)"
$ws.Range("D6").Value = "(1.02)"

$ws.Range("E4").Value = "(0.39)"
$ws.Range("E6").Value = "(0.07)"

$ws.Range("F4").Value = "(1.13)"
$ws.Range("F6").Value = "(1.53)"

$ws.Range("G4").Value = "(2.12)"
$ws.Range("G6").Value = "(1.74)"

$ws.Range("H4").Value = "(2.34)"
$ws.Range("H6").Value = "(1.09)"

$ws.Range("I4").Value = "(1.88)"
$ws.Range("I6").Value = "(2.4)"

$ws.Range("J4").Value = "(0.05)"
$ws.Range("J6").Value = "(0.13)"
